$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("2025-07-07 15:48:43", "Unknown", "Unknown"),
    @("2025-07-07 15:49:02", "Unknown", "Unknown"),
    @("2025-07-07 15:49:44", "Unknown", "Unknown"),
    @("2025-07-07 15:50:06", "xuan_3", "Success")
)

$row = 7
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row++
}
